$d = $word.ActiveDocument

# Each "<id>...</id>" paragraph in this document is split across three
# runs: "<id>" (Courier New / color 7f6000 / sz 18), the bare id text
# (plain formatting), and "</id>" (Courier New / color 7f6000 / sz 18
# again). The edit merges each of those three runs into a single run
# (keeping the "<id>" run's formatting) containing the whole
# "<id>p064r_N</id>" string - renumbering sequentially (occurrence #5 and
# #8 were missing the "r" in the original id and get corrected too).
$counter = 0
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $t = $r.Text
    if ($t.StartsWith("<id>") -and $t.Contains("</id>")) {
        $counter = $counter + 1
        $newText = "<id>p064r_" + $counter + "</id>"

        # Limit the range to just the "<id>...</id>" text (exclude the
        # trailing paragraph mark / anything after </id>) so formatting of
        # the paragraph-mark run is left untouched.
        $closeIdx = $t.IndexOf("</id>") + 5
        $idRange = $d.Range($r.Start, $r.Start + $closeIdx)

        # Setting .Text collapses the whole range into a single run that
        # carries the formatting of the range's first character (the
        # "<id>" run: Courier New, color 7f6000, sz 18). When the target
        # text happens to be byte-identical to the current text (true for
        # the occurrences that were already correctly numbered) a direct
        # assignment is a no-op and the three runs are left unmerged, so
        # first stage the range through a throwaway one-character value to
        # force the rebuild, then set the real text on the now-single
        # character run.
        $idRange.Text = "X"
        $stagedRange = $d.Range($r.Start, $r.Start + 1)
        $stagedRange.Text = $newText
    }
}
